$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 158 (shifts rows 158:249 down to 159:250,
# and grows the used range from A1:R249 to A1:R250).
$ws.Rows.Item(158).Insert()

# Populate the newly inserted row 158 with the new record.
$ws.Range("A158").Value = 10
$ws.Range("B158").Value = "Vega Modelo de Temuco"
$ws.Range("C158").Value = "La Araucanía"
$ws.Range("D158").Value = 44488
$ws.Range("E158").Value = 9
$ws.Range("F158").Value = 100112008
$ws.Range("G158").Value = "Coliflor"
$ws.Range("H158").Value = "Sin especificar"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 2850
$ws.Range("K158").Value = 800
$ws.Range("L158").Value = 800
$ws.Range("M158").Value = 800
$ws.Range("N158").Value = "$/unidad"
$ws.Range("O158").Value = "Región Metropolitana"
$ws.Range("P158").Value = 800
$ws.Range("Q158").Value = 1
$ws.Range("R158").Value = "Hortaliza"
